$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2:G20").NumberFormat = "@"

$ws.Range("G1").Value = "product_code"

$ws.Range("G15").Value = "6-41"
$ws.Range("G16").Value = "6-42"
$ws.Range("G17").Value = "6-43"
$ws.Range("G5").Value = "885-2"
$ws.Range("G6").Value = "885-5"

$ws.Range("B10").Select() | Out-Null
